$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update dependency version text (ruby/mb version bump) ---
$ws.Range("B1").Value = "EOL mb version"
$ws.Range("B2").Value = "v1.7"
$ws.Range("B3").Value = "v1.7"
$ws.Range("C6").Value = "7.1.0"
$ws.Range("C7").Value = "v6.9.1"

# --- Toggle bold on the "mb version" column ---
# Header becomes bold, the rest of the column loses its bold formatting.
$ws.Range("B1").Font.Bold = $true
$ws.Range("B2").Font.Bold = $false
$ws.Range("B3").Font.Bold = $false
$ws.Range("B4").Font.Bold = $false
$ws.Range("B5").Font.Bold = $false
$ws.Range("B6").Font.Bold = $false
$ws.Range("B7").Font.Bold = $false
$ws.Range("B8").Font.Bold = $false
$ws.Range("B9").Font.Bold = $false

# --- Resize column B to fit the new (shorter) content ---
$ws.Columns("B").ColumnWidth = 13.5
